$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D and E columns hold text-formatted values (e.g. "27.390.84", "  -2.87%  ").
# Set NumberFormat to text ("@") first so Excel does not auto-convert these
# numeric-looking strings into actual numbers / percentages.
$priceRange = $ws.Range("D2:D51")
$volRange = $ws.Range("E2:E51")
$priceRange.NumberFormat = "@"
$volRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.390.84'
$ws.Range("E2").Value = '  -2.87%  '
$ws.Range("D3").Value = '1.738.88'
$ws.Range("E3").Value = '  -3.57%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '322.90'
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4242'
$ws.Range("E7").Value = '  -9.08%  '
$ws.Range("D8").Value = '0.3610'
$ws.Range("E8").Value = '  -2.72%  '
$ws.Range("D9").Value = '45.41'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '0.07426'
$ws.Range("E10").Value = '  -3.24%  '
$ws.Range("D11").Value = '1.114'
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("D12").Value = '0.9997'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '21.56'
$ws.Range("E13").Value = '  -4.73%  '
$ws.Range("D14").Value = '6.071'
$ws.Range("E14").Value = '  -4.81%  '
$ws.Range("D15").Value = '7.169'
$ws.Range("E15").Value = '  -3.16%  '
$ws.Range("D16").Value = '1.729.35'
$ws.Range("E16").Value = '  -3.63%  '
$ws.Range("D17").Value = '0.00001060'
$ws.Range("E17").Value = '  -3.32%  '
$ws.Range("D18").Value = '87.33'
$ws.Range("E18").Value = '  +5.60%  '
$ws.Range("D19").Value = '0.06011'
$ws.Range("E19").Value = '  -10.62%  '
$ws.Range("D20").Value = '0.9994'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '16.82'
$ws.Range("E21").Value = '  -3.72%  '
$ws.Range("D22").Value = '6.084'
$ws.Range("E22").Value = '  -5.37%  '
$ws.Range("D23").Value = '0.5236'
$ws.Range("E23").Value = '  -3.99%  '
$ws.Range("D24").Value = '27.395.54'
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("D25").Value = '11.33'
$ws.Range("E25").Value = '  -5.03%  '
$ws.Range("D26").Value = '2.385'
$ws.Range("E26").Value = '  -1.34%  '
$ws.Range("D27").Value = '20.18'
$ws.Range("E27").Value = '  -3.29%  '
$ws.Range("D28").Value = '2.374'
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("D29").Value = '149.15'
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("D30").Value = '1.927.49'
$ws.Range("E30").Value = '  -3.85%  '
$ws.Range("D31").Value = '126.47'
$ws.Range("E31").Value = '  -5.93%  '
$ws.Range("D32").Value = '1.202'
$ws.Range("E32").Value = '  -5.28%  '
$ws.Range("D33").Value = '5.649'
$ws.Range("E33").Value = '  -4.56%  '
$ws.Range("D34").Value = '0.09081'
$ws.Range("E34").Value = '  -5.94%  '
$ws.Range("D35").Value = '3.660'
$ws.Range("E35").Value = '  -9.51%  '
$ws.Range("D36").Value = '12.89'
$ws.Range("E36").Value = '  +5.25%  '
$ws.Range("D37").Value = '0.2138'
$ws.Range("E37").Value = '  -4.85%  '
$ws.Range("D38").Value = '5.057'
$ws.Range("E38").Value = '  -4.12%  '
$ws.Range("D39").Value = '0.02251'
$ws.Range("E39").Value = '  -5.30%  '
$ws.Range("D40").Value = '0.06056'
$ws.Range("E40").Value = '  -5.43%  '
$ws.Range("D41").Value = '0.6384'
$ws.Range("E41").Value = '  -5.15%  '
$ws.Range("D42").Value = '1.188'
$ws.Range("E42").Value = '  -4.14%  '
$ws.Range("D43").Value = '7.937'
$ws.Range("E43").Value = '  -2.76%  '
$ws.Range("D44").Value = '0.9991'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '1.408'
$ws.Range("E45").Value = '  -7.78%  '
$ws.Range("D46").Value = '13.56'
$ws.Range("E46").Value = '  -4.39%  '
$ws.Range("D47").Value = '3.723'
$ws.Range("E47").Value = '  -3.10%  '
$ws.Range("D48").Value = '0.5830'
$ws.Range("E48").Value = '  -5.72%  '
$ws.Range("D49").Value = '124.77'
$ws.Range("E49").Value = '  -4.11%  '
$ws.Range("D50").Value = '1.957'
$ws.Range("E50").Value = '  -5.45%  '
$ws.Range("D51").Value = '0.06843'
$ws.Range("E51").Value = '  -4.19%  '
